$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2089.7073
$ws.Range("J17").Value = 2089.7073
$ws.Range("L17").Value = 6269.1219
$ws.Range("N17").Value = -6605.1219
$ws.Range("H19").Value = 1069.2727
$ws.Range("I19").Value = 999.5
$ws.Range("J19").Value = 1255.3334
$ws.Range("K19").Value = 999.5
$ws.Range("L19").Value = 1255.3334
$ws.Range("M19").Value = -824.5
$ws.Range("N19").Value = -1605.3334
$ws.Range("H51").Value = 9631.786
$ws.Range("I51").Value = 7333.3335
$ws.Range("J51").Value = 10258.637
$ws.Range("K51").Value = 7333.3335
$ws.Range("L51").Value = 10258.637
$ws.Range("M51").Value = -6849.3335
$ws.Range("N51").Value = -11226.637
$ws.Range("H74").Value = 8582.091
$ws.Range("I74").Value = 8425.375
$ws.Range("K74").Value = 8425.375
$ws.Range("M74").Value = -7489.375
$ws.Range("H76").Value = 6736.1816
$ws.Range("I76").Value = 5399.6665
$ws.Range("J76").Value = 7237.375
$ws.Range("K76").Value = 5399.6665
$ws.Range("L76").Value = 7237.375
$ws.Range("M76").Value = -5084.6665
$ws.Range("N76").Value = -7867.375
$ws.Range("H77").Value = 8582.091
$ws.Range("I77").Value = 8425.375
$ws.Range("K77").Value = 42126.875
$ws.Range("M77").Value = -37446.875
$ws.Range("H79").Value = 6736.1816
$ws.Range("I79").Value = 5399.6665
$ws.Range("J79").Value = 7237.375
$ws.Range("K79").Value = 5399.6665
$ws.Range("L79").Value = 7237.375
$ws.Range("M79").Value = -4307.6665
$ws.Range("N79").Value = -9421.375
$ws.Range("H80").Value = 882.8889
$ws.Range("J80").Value = 1124.5
$ws.Range("L80").Value = 3373.5
$ws.Range("N80").Value = -5369.5
$ws.Range("H83").Value = 882.8889
$ws.Range("J83").Value = 1124.5
$ws.Range("L83").Value = 10120.5
$ws.Range("N83").Value = -20104.5
$ws.Range("H86").Value = 6541.324
$ws.Range("I86").Value = 6127.0415
$ws.Range("J86").Value = 7306.154
$ws.Range("K86").Value = 6127.0415
$ws.Range("L86").Value = 7306.154
$ws.Range("M86").Value = -5004.0415
$ws.Range("N86").Value = -9552.154
$ws.Range("H89").Value = 6541.324
$ws.Range("I89").Value = 6127.0415
$ws.Range("J89").Value = 7306.154
$ws.Range("K89").Value = 30635.2075
$ws.Range("L89").Value = 36530.77
$ws.Range("M89").Value = -25019.2075
$ws.Range("N89").Value = -47762.77
$ws.Range("H138").Value = 3248.7368
$ws.Range("J138").Value = 3530.625
$ws.Range("L138").Value = 10591.875
$ws.Range("N138").Value = -20871.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7528.5747
$ws.Range("I32").Value = 2829.8809
$ws.Range("K32").Value = 2829.8809
$ws.Range("M32").Value = -2542.8809
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4586
$ws.Range("H61").Value = 4784.9165
$ws.Range("I61").Value = 4058.3333
$ws.Range("J61").Value = 5511.5
$ws.Range("K61").Value = 4058.3333
$ws.Range("L61").Value = 5511.5
$ws.Range("M61").Value = -3846.3333
$ws.Range("N61").Value = -5935.5
$ws.Range("H132").Value = 7308.6816
$ws.Range("I132").Value = 2199.6428
$ws.Range("J132").Value = 16249.5
$ws.Range("K132").Value = 6598.928400000001
$ws.Range("L132").Value = 48748.5
$ws.Range("M132").Value = -4068.928400000001
$ws.Range("N132").Value = -53808.5
$ws.Range("H136").Value = 4784.9165
$ws.Range("I136").Value = 4058.3333
$ws.Range("J136").Value = 5511.5
$ws.Range("K136").Value = 12174.9999
$ws.Range("L136").Value = 16534.5
$ws.Range("M136").Value = -9624.999899999999
$ws.Range("N136").Value = -21634.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 713.3415
$ws.Range("I94").Value = 706.7568
$ws.Range("K94").Value = 706.7568
$ws.Range("M94").Value = -255.7568

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 379.63635
$ws.Range("I7").Value = 273.88235
$ws.Range("J7").Value = 492
$ws.Range("K7").Value = 273.88235
$ws.Range("L7").Value = 492
$ws.Range("M7").Value = -160.88235
$ws.Range("N7").Value = -718
$ws.Range("H16").Value = 962.0769
$ws.Range("I16").Value = 522.2727
$ws.Range("K16").Value = 522.2727
$ws.Range("M16").Value = -235.2727
$ws.Range("H31").Value = 9446.416999999999
$ws.Range("I31").Value = 3539.7222
$ws.Range("K31").Value = 3539.7222
$ws.Range("M31").Value = -3244.7222
$ws.Range("H34").Value = 9446.416999999999
$ws.Range("I34").Value = 3539.7222
$ws.Range("K34").Value = 3539.7222
$ws.Range("M34").Value = -3337.7222
$ws.Range("H41").Value = 12424.75
$ws.Range("I41").Value = 5274.75
$ws.Range("J41").Value = 33874.75
$ws.Range("K41").Value = 5274.75
$ws.Range("L41").Value = 33874.75
$ws.Range("M41").Value = -4846.75
$ws.Range("N41").Value = -34730.75
$ws.Range("H42").Value = 4778
$ws.Range("I42").Value = 5056
$ws.Range("J42").Value = 4500
$ws.Range("K42").Value = 5056
$ws.Range("L42").Value = 4500
$ws.Range("M42").Value = -4463
$ws.Range("N42").Value = -5686
$ws.Range("H86").Value = 7979
$ws.Range("I86").Value = 7931.6665
$ws.Range("J86").Value = 8050
$ws.Range("K86").Value = 7931.6665
$ws.Range("L86").Value = 8050
$ws.Range("M86").Value = -6808.6665
$ws.Range("N86").Value = -10296
$ws.Range("H89").Value = 7979
$ws.Range("I89").Value = 7931.6665
$ws.Range("J89").Value = 8050
$ws.Range("K89").Value = 39658.3325
$ws.Range("L89").Value = 40250
$ws.Range("M89").Value = -34042.3325
$ws.Range("N89").Value = -51482
$ws.Range("H113").Value = 962.0769
$ws.Range("I113").Value = 522.2727
$ws.Range("K113").Value = 522.2727
$ws.Range("M113").Value = 1647.7273
$ws.Range("H122").Value = 3249.6155
$ws.Range("I122").Value = 3175.8572
$ws.Range("J122").Value = 3559.4
$ws.Range("K122").Value = 9527.571599999999
$ws.Range("L122").Value = 10678.2
$ws.Range("M122").Value = -7077.571599999999
$ws.Range("N122").Value = -15578.2
$ws.Range("H132").Value = 3782.6155
$ws.Range("I132").Value = 3379.9092
$ws.Range("J132").Value = 5997.5
$ws.Range("K132").Value = 10139.7276
$ws.Range("L132").Value = 17992.5
$ws.Range("M132").Value = -7609.7276
$ws.Range("N132").Value = -23052.5
$ws.Range("H134").Value = 6054.3887
$ws.Range("I134").Value = 4927.2856
$ws.Range("J134").Value = 9999.25
$ws.Range("K134").Value = 14781.8568
$ws.Range("L134").Value = 29997.75
$ws.Range("M134").Value = -12246.8568
$ws.Range("N134").Value = -35067.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 20.333334
$ws.Range("I7").Value = 17.666666
$ws.Range("K7").Value = 52.999998
$ws.Range("M7").Value = 59.000002
$ws.Range("H55").Value = 3280.4
$ws.Range("I55").Value = 801.3333
$ws.Range("J55").Value = 6999
$ws.Range("K55").Value = 2403.9999
$ws.Range("L55").Value = 20997
$ws.Range("M55").Value = -2226.9999
$ws.Range("N55").Value = -21351
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H131").Value = 651952.5
$ws.Range("J131").Value = 1403068.9
$ws.Range("L131").Value = 4209206.699999999
$ws.Range("N131").Value = -4219286.699999999
$ws.Range("H140").Value = 1853.6285
$ws.Range("I140").Value = 1520.4546
$ws.Range("J140").Value = 2006.3334
$ws.Range("K140").Value = 4561.3638
$ws.Range("L140").Value = 6019.0002
$ws.Range("M140").Value = 618.6361999999999
$ws.Range("N140").Value = -16379.0002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8294.52
$ws.Range("J70").Value = 8559.380999999999
$ws.Range("L70").Value = 8559.380999999999
$ws.Range("N70").Value = -9099.380999999999
$ws.Range("H73").Value = 8294.52
$ws.Range("J73").Value = 8559.380999999999
$ws.Range("L73").Value = 8559.380999999999
$ws.Range("N73").Value = -10431.381
$ws.Range("H97").Value = 581
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 2000
$ws.Range("N97").Value = -2992
$ws.Range("H132").Value = 5991.5
$ws.Range("I132").Value = 2874.4285
$ws.Range("J132").Value = 9108.571
$ws.Range("K132").Value = 8623.2855
$ws.Range("L132").Value = 27325.713
$ws.Range("M132").Value = -6093.2855
$ws.Range("N132").Value = -32385.713
$ws.Range("H135").Value = 98992.625
$ws.Range("J135").Value = 98992.625
$ws.Range("L135").Value = 98992.625
$ws.Range("N135").Value = -109132.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1781.0294
$ws.Range("I16").Value = 1807.875
$ws.Range("J16").Value = 1716.6
$ws.Range("K16").Value = 1807.875
$ws.Range("L16").Value = 1716.6
$ws.Range("M16").Value = -1637.875
$ws.Range("N16").Value = -2056.6
$ws.Range("H46").Value = 2826.5908
$ws.Range("J46").Value = 5468.5
$ws.Range("L46").Value = 5468.5
$ws.Range("N46").Value = -5844.5
$ws.Range("H93").Value = 2543.9092
$ws.Range("I93").Value = 2653.3
$ws.Range("K93").Value = 2653.3
$ws.Range("M93").Value = -1405.3

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H100").Value = 828.5833
$ws.Range("J100").Value = 924.625
$ws.Range("L100").Value = 1849.25
$ws.Range("N100").Value = -2931.25
$ws.Range("H132").Value = 3953.8
$ws.Range("I132").Value = 2948.6667
$ws.Range("K132").Value = 8846.000100000001
$ws.Range("M132").Value = -6316.000100000001
$ws.Range("H136").Value = 4301.3076
$ws.Range("I136").Value = 4326.4165
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 12979.2495
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -10429.2495
$ws.Range("N136").Value = -17100
